# Source data extraction from sources 1 through 9
# Fill in the consensus-mechanism comparison table (rows 2-10, cols B-G)
# on the "Data" sheet, and tidy up the now-unused "Percent" cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Row 2: Proof of Work ---------------------------------------------
$ws.Range("B2").Formula = "=AVERAGE(3.78,7,56,4,100, 60, 56)"
$ws.Range("C2:G2").Value = "N/A"

# --- Row 3: Proof of Stake ---------------------------------------------
$ws.Range("B3").Formula = "=AVERAGE(40.56,30,250,1000)"
$ws.Range("C3:G3").Value = "N/A"

# --- Row 4: Delegated Proof of Stake -----------------------------------
$ws.Range("B4:G4").Value = "N/A"

# --- Row 5: Proof of History --------------------------------------------
$ws.Range("B5").Value = 65000
$ws.Range("C5:G5").Value = "N/A"

# --- Row 6: Proof of Stake with Byzantine Fault Tolerance ---------------
$ws.Range("B6:G6").Value = "N/A"

# --- Row 7: Proof of History with Proof of Stake ------------------------
$ws.Range("B7").Value = 65000
$ws.Range("C7:G7").Value = "N/A"

# --- Row 8: zk-proof ------------------------------------------------------
$ws.Range("B8:G8").Value = "N/A"

# --- Row 9: Sharding --------------------------------------------------------
$ws.Range("B9:G9").Value = "N/A"

# --- Row 10: DAGs -------------------------------------------------------------
$ws.Range("B10:G10").Value = "N/A"

# The "% of nodes..." column (E) previously carried the unused "Percent"
# cell style (blank cells formatted as percent) - now that every row has
# real data, drop that now-orphaned style from the workbook.
$ws.Range("E2:E10").Style = "Normal"
$wb.Styles.Item("Percent").Delete()

# Move the active selection as left by the editor.
$ws.Range("D16").Select() | Out-Null
